# Updated symbol list (coin prices / rankings refresh).
# D-column prices are stored as literal text in the workbook (General
# format), so new values are written with a leading apostrophe
# (Excel quote-prefix) to keep them as text and preserve exact
# formatting (trailing zeros, leading zeros, no scientific notation)
# instead of letting Excel auto-convert them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.27"
$ws.Range("D3").Value = "'23.97"
$ws.Range("D4").Value = "'5.268"
$ws.Range("D5").Value = "'0.05823"
$ws.Range("D6").Value = "'6.466"
$ws.Range("D7").Value = "'3.231"
$ws.Range("D8").Value = "'0.8081"
$ws.Range("D9").Value = "'0.8824"
$ws.Range("D10").Value = "'0.1388"
$ws.Range("D11").Value = "'0.07130"
$ws.Range("D12").Value = "'0.03079"
$ws.Range("D13").Value = "'0.03042"
$ws.Range("D14").Value = "'0.09334"
$ws.Range("D15").Value = "'3.813"
$ws.Range("D16").Value = "'0.001536"
$ws.Range("D17").Value = "'0.04703"
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = "'0.006208"
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").Value = "'0.001259"
$ws.Range("E19").Value = '18BitKanKAN'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").Value = "'0.004077"
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Value = "'0.00008701"
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").Value = "'3.540"
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").Value = "'2.162"
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").Value = "'0.01035"
$ws.Range("E24").Value = '23OneONEBestin24h'
$ws.Range("D26").Value = "'0.1315"
$ws.Range("D28").Value = "'0.0002329"
$ws.Range("D40").Value = "'0.03845"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.006280"
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1054"
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.002541"
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = "'0.007238"
$ws.Range("D45").Value = "'0.00005327"
$ws.Range("D47").Value = "'0.5350"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'
$ws.Range("D48").Value = "'0.003868"
$ws.Range("E48").Value = '47BOLOBOLO'
